# ---------------------------------------------------------------------------
# Applies the "2nd revision board" documentation update to resources.xlsx:
#   - DMA sheet: fill in DMA1 / DMA2 request-mapping tables
#   - Interrupt Priorities sheet: add EXTI_12 (shares IRQ w/ EXTI_11), a new
#     "Priority 1" section for the IMU's SPI3 RX/TX, and a new "Priority 2"
#     section for the Start/Pause button (EXTI_6)
#   - Make "Interrupt Priorities" the active/selected sheet
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Excel COM alignment constant values (xlHAlign*, xlVAlign*)
$xlHAlignCenter = -4108
$xlHAlignRight  = -4152
$xlVAlignCenter = -4108

# ===========================================================================
# DMA sheet
# ===========================================================================
$dma = $wb.Worksheets.Item("DMA")

# --- DMA1 table -------------------------------------------------------------
$dma.Cells.Item(1,1).Value = "DMA1"
$dma.Range("A1").Style = "Heading 1"
$dma.Range("A1").HorizontalAlignment = $xlHAlignRight
$dma.Rows.Item(1).RowHeight = 20.4
$dma.Rows.Item(2).RowHeight = 15

$dma.Cells.Item(3,1).Value = "Request"
$dma.Range("A3").Style = "Good"
$dma.Range("A3").HorizontalAlignment = $xlHAlignRight

$dma.Cells.Item(3,2).Value = "Channel 1"
$dma.Cells.Item(3,3).Value = "Channel 2"
$dma.Cells.Item(3,4).Value = "Channel 3"
$dma.Cells.Item(3,5).Value = "Channel 4"
$dma.Cells.Item(3,6).Value = "Channel 5"
$dma.Cells.Item(3,7).Value = "Channel 6"
$dma.Cells.Item(3,8).Value = "Channel 7"
$dma.Range("B3:H3").Style = "Good"

for ($r = 4; $r -le 11; $r++) {
    $dma.Cells.Item($r,1).Value = $r - 4
}

$dma.Cells.Item(4,2).Value = "ADC1"
$dma.Cells.Item(4,3).Value = "ADC2"
$dma.Cells.Item(4,4).Value = "ADC3"

$dma.Cells.Item(6,5).Value = "USART1_TX"
$dma.Cells.Item(6,6).Value = "USART1_RX"
$dma.Cells.Item(6,7).Value = "USART2_RX"
$dma.Cells.Item(6,8).Value = "USART2_TX"

$rngDma1 = $dma.Range("A4:H11")
$rngDma1.HorizontalAlignment = $xlHAlignCenter
$rngDma1.VerticalAlignment = $xlVAlignCenter

# --- DMA2 table --------------------------------------------------------------
$dma.Cells.Item(13,1).Value = "DMA2"
$dma.Range("A13").Style = "Heading 1"
$dma.Range("A13").HorizontalAlignment = $xlHAlignRight
$dma.Rows.Item(13).RowHeight = 20.4
$dma.Rows.Item(14).RowHeight = 15

$dma.Cells.Item(15,1).Value = "Request"
$dma.Range("A15").Style = "Good"
$dma.Range("A15").HorizontalAlignment = $xlHAlignRight

$dma.Cells.Item(15,2).Value = "Channel 1"
$dma.Cells.Item(15,3).Value = "Channel 2"
$dma.Cells.Item(15,4).Value = "Channel 3"
$dma.Cells.Item(15,5).Value = "Channel 4"
$dma.Cells.Item(15,6).Value = "Channel 5"
$dma.Cells.Item(15,7).Value = "Channel 6"
$dma.Cells.Item(15,8).Value = "Channel 7"
$dma.Range("B15:H15").Style = "Good"

for ($r = 16; $r -le 23; $r++) {
    $dma.Cells.Item($r,1).Value = $r - 16
}

$dma.Cells.Item(19,2).Value = "SPI3_RX"
$dma.Cells.Item(19,3).Value = "SPI3_TX"

$dma.Cells.Item(21,7).Value = "I2C1_RX"
$dma.Cells.Item(21,8).Value = "I2C1_TX"

$dma.Cells.Item(23,5).Value = "SDMMC1"

$rngDma2 = $dma.Range("A17:H23")
$rngDma2.HorizontalAlignment = $xlHAlignCenter
$rngDma2.VerticalAlignment = $xlVAlignCenter

$dma.Range("A16").HorizontalAlignment = $xlHAlignCenter
$dma.Range("A16").VerticalAlignment = $xlVAlignCenter

$dma.Range("D16").HorizontalAlignment = $xlHAlignCenter
$dma.Range("D16").VerticalAlignment = $xlVAlignCenter
$dma.Range("F16").HorizontalAlignment = $xlHAlignCenter
$dma.Range("F16").VerticalAlignment = $xlVAlignCenter

# --- column widths (character units, offset to match the engine's internal
#     pixel-rounding of the ColumnWidth property) ---------------------------
$dma.Columns.Item(1).ColumnWidth = 7.721354166666667
$dma.Range("B1:D1").ColumnWidth = 8.166666666666666
$dma.Columns.Item(5).ColumnWidth = 9.608072916666666
$dma.Range("F1:G1").ColumnWidth = 9.721354166666666
$dma.Columns.Item(8).ColumnWidth = 9.608072916666666

$dma.Range("E16").Select()

# ===========================================================================
# Interrupt Priorities sheet
# ===========================================================================
$ip = $wb.Worksheets.Item("Interrupt Priorities")

# Existing "Priority 0" section gets a new row: EXTI_12
$ip.Cells.Item(5,1).Value = "EXTI_12"
$ip.Cells.Item(5,2).Value = "Flush SD (will be moved later to EXTI_9)"
$ip.Cells.Item(5,3).Value = "Shares the same IRQ as EXTI_11"

# New "Priority 1" section
$ip.Cells.Item(7,1).Value = "Priority 1"
$ip.Range("A7").Style = "Heading 1"
$ip.Rows.Item(7).RowHeight = 20.4
$ip.Rows.Item(8).RowHeight = 15

$ip.Cells.Item(9,1).Value = "SPI3 RX"
$ip.Cells.Item(9,2).Value = "IMU"

$ip.Cells.Item(10,1).Value = "SPI3 TX"
$ip.Cells.Item(10,2).Value = "IMU"

# New "Priority 2" section
$ip.Cells.Item(12,1).Value = "Priority 2"
$ip.Range("A12").Style = "Heading 1"
$ip.Rows.Item(12).RowHeight = 20.4
$ip.Rows.Item(13).RowHeight = 15

$ip.Cells.Item(14,1).Value = "EXTI_6"
$ip.Cells.Item(14,2).Value = "Start/Pause button"

$ip.Columns.Item(2).ColumnWidth = 32.60807291666667
$ip.Columns.Item(3).ColumnWidth = 26.498697916666664

$ip.Activate()
$ip.Range("B10").Select()

Write-Host "edit applied"
